# CreatedOn kolonu eklendi, simdilik geri bildirim sonucu hareket edecegiz
#
# The "Survey" table (columns E:F, headed at E1/F1) gets a new
# "createdOn     datetime" row inserted right after "score float(5)"
# (row 5) and before the existing "IsActive     bit" row, which shifts
# down by one row (from E6 to E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current E6 value ("IsActive     bit") down into E7 ...
$ws.Range("E7").Value = $ws.Range("E6").Value2

# ... and replace E6 with the new "createdOn" column definition.
$ws.Range("E6").Value = "createdOn     datetime"

# Update the active window selection to match the authored view state.
# (Note: this engine's Window.ScrollColumn/ScrollRow — and freezing/
# unfreezing panes — do not persist a bare `topLeftCell` on <sheetView>;
# only the active cell/selection is settable here.)
$ws.Range("E8").Select()
